# autoIt and trace log
# The "Google" sheet test-automation log is updated after a re-run of the
# AutoIt search-automation / trace logging:
#   - Test case 1 (row 3, "Zsinj" search) no longer has a captured
#     First Result / URL pair - clear E3 and F3.
#   - Test case 4 (row 6, "windows" search) now reports a failure instead
#     of a pass - update D6's Status from "pass" to "fail".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Google")

$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("D6").Value = "fail"
